# Set up run modes for test cases:
# 1) Insert a new "test_suite" sheet at the front describing which test
#    cases should run and in what mode.
# 2) Add two more customer rows (vicky/gaurav) to AddCustomerTest.
# 3) Fix the currency label on OpenAccountTest ("rupee" -> "Rupee").

$wb = $excel.ActiveWorkbook

# --- 1. new "test_suite" sheet, inserted before the first existing sheet ---
$testSuite = $wb.Worksheets.Add($wb.Worksheets.Item(1))
$testSuite.Name = "test_suite"

$testSuite.Range("A1").Value = "TCID"
$testSuite.Range("B1").Value = "Runmode"
$testSuite.Range("A2").Value = "BankManagerLoginTest"
$testSuite.Range("B2").Value = "Y"
$testSuite.Range("A3").Value = "AddCustomerTest"
$testSuite.Range("B3").Value = "Y"
$testSuite.Range("A4").Value = "OpenAccountTest"
$testSuite.Range("B4").Value = "N"

$testSuite.Columns.Item(1).ColumnWidth = 23.42
$testSuite.Columns.Item(2).ColumnWidth = 14.92

$testSuite.Range("B1:B4").Select()

# --- 2. add rows to AddCustomerTest (re-fetch by name after sheet insert) ---
$addCustomer = $wb.Worksheets.Item("AddCustomerTest")

$addCustomer.Range("A4").Value = "vicky"
$addCustomer.Range("B4").Value = "thopate"
$addCustomer.Range("C4").Value = 412303
$addCustomer.Range("D4").Value = "Customer added successfully"

$addCustomer.Range("A5").Value = "gaurav"
$addCustomer.Range("B5").Value = "yadav"
$addCustomer.Range("C5").Value = 412301
$addCustomer.Range("D5").Value = "Customer added successfully"

# --- 3. fix currency label on OpenAccountTest ---
$openAccount = $wb.Worksheets.Item("OpenAccountTest")
$openAccount.Range("B2").Value = "Rupee"
$openAccount.Range("B2").Select()

# --- final selection / active sheet state ---
$addCustomer = $wb.Worksheets.Item("AddCustomerTest")
$addCustomer.Activate()
$addCustomer.Range("E1:E7").Select()
